# Update crypto price/volume values per latest scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'58.800.70"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  -0.25%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'2.308.75"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'538.46"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -2.01%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'132.39"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +0.33%  "
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +0.04%  "
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +2.35%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'2.307.86"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -0.32%  "
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -1.87%  "
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -1.30%  "
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +0.90%  "
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -0.61%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'23.73"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -1.43%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'2.720.47"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -0.28%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'58.661.54"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -0.34%  "
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  -0.44%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'2.291.40"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'10.64"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -0.82%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'4.18"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -3.52%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'313.70"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -0.64%  "
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +1.91%  "
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -0.01%  "
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -0.99%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'0.172"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +0.76%  "
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +0.12%  "
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -2.14%  "
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -2.33%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'171.32"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  +0.98%  "
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -2.22%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'0.0₃0734"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +0.42%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'1.14"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  +2.16%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'5.86"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +0.79%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'0.386"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +0.29%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'17.92"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +0.54%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'1.30"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +2.65%  "
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -0.01%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'4.06"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +1.33%  "
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +0.17%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'295.74"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -2.77%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'141.10"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -0.35%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'3.45"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +0.03%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.0961"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +0.94%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.0496"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -1.31%  "
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -0.73%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'18.34"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -1.91%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'0.0211"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -2.26%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'11.00"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -0.24%  "
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +0.06%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'1.51"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  +0.48%  "
$c.Style = "Normal"
